$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.078.05"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.057.55"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'249.47"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "'0.672"
$ws.Range("E6").Value = "  +1.84%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'56.06"
$ws.Range("E8").Value = "  +9.03%  "
$ws.Range("D9").Value = "'60.89"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").Value = "'0.378"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").Value = "'0.0804"
$ws.Range("E11").Value = "  +7.98%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "'15.22"
$ws.Range("E13").Value = "  +4.60%  "
$ws.Range("D14").Value = "2.361.48"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "'0.814"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").Value = "'5.33"
$ws.Range("E16").Value = "  +3.51%  "
$ws.Range("D17").Value = "2.060.64"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").Value = "37.001.09"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "'74.41"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").Value = "0.0₃0920"
$ws.Range("E20").Value = "  +11.50%  "
$ws.Range("D21").Value = "'14.27"
$ws.Range("E21").Value = "  +7.06%  "
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("D23").Value = "'237.78"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("D26").Value = "'171.69"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").Value = "'9.08"
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("D28").Value = "'20.14"
$ws.Range("E28").Value = "  -4.62%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  +1.70%  "
$ws.Range("D31").Value = "'4.60"
$ws.Range("E31").Value = "  +2.09%  "
$ws.Range("D32").Value = "'1.07"
$ws.Range("E32").Value = "  -6.86%  "
$ws.Range("D33").Value = "'0.0628"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("D34").Value = "'4.40"
$ws.Range("E34").Value = "  +6.83%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'0.0869"
$ws.Range("E36").Value = "  -6.05%  "
$ws.Range("E37").Value = "  -2.86%  "
$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").Value = "'0.106"
$ws.Range("E40").Value = "  +22.62%  "
$ws.Range("D41").Value = "'4.66"
$ws.Range("E41").Value = "  +57.28%  "
$ws.Range("D42").Value = "'18.03"
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").Value = "'97.03"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  +11.34%  "
$ws.Range("D48").Value = "'2.45"
$ws.Range("E48").Value = "  +7.87%  "
$ws.Range("D49").Value = "1.300.01"
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("D50").Value = "'2.91"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "'6.86"
$ws.Range("E51").Value = "  -2.27%  "
